$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 526.4211
$ws.Range("I39").Value = 107.28571
$ws.Range("J39").Value = 1700
$ws.Range("K39").Value = 321.85713
$ws.Range("L39").Value = 5100
$ws.Range("M39").Value = -25.85712999999998
$ws.Range("N39").Value = -5692

$ws.Range("H42").Value = 104.625
$ws.Range("I42").Value = 91
$ws.Range("K42").Value = 273
$ws.Range("M42").Value = -43

$ws.Range("H132").Value = 553048.1
$ws.Range("I132").Value = 2288.0132
$ws.Range("J132").Value = 3772876.5
$ws.Range("K132").Value = 6864.0396
$ws.Range("L132").Value = 11318629.5
$ws.Range("M132").Value = -4334.0396
$ws.Range("N132").Value = -11323689.5

$ws.Range("H135").Value = 35172.266
$ws.Range("I135").Value = 40179.54
$ws.Range("J135").Value = 2625
$ws.Range("K135").Value = 361615.86
$ws.Range("L135").Value = 23625
$ws.Range("M135").Value = -359080.86
$ws.Range("N135").Value = -28695

$ws.Range("H137").Value = 3850970.8
$ws.Range("I137").Value = 5559477
$ws.Range("J137").Value = 6831.75
$ws.Range("K137").Value = 16678431
$ws.Range("L137").Value = 20495.25
$ws.Range("M137").Value = -16675881
$ws.Range("N137").Value = -25595.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2791
$ws.Range("I2").Value = 2682.8333
$ws.Range("K2").Value = 2682.8333
$ws.Range("M2").Value = -2569.8333

$ws.Range("H32").Value = 11530.57
$ws.Range("I32").Value = 8983.468999999999
$ws.Range("J32").Value = 22389.264
$ws.Range("K32").Value = 8983.468999999999
$ws.Range("L32").Value = 22389.264
$ws.Range("M32").Value = -8696.468999999999
$ws.Range("N32").Value = -22963.264

$ws.Range("H45").Value = 1626.625
$ws.Range("I45").Value = 1395.2307
$ws.Range("J45").Value = 2629.3333
$ws.Range("K45").Value = 1395.2307
$ws.Range("L45").Value = 2629.3333
$ws.Range("M45").Value = -1018.2307
$ws.Range("N45").Value = -3383.3333

$ws.Range("H74").Value = 8828259
$ws.Range("I74").Value = 13376108
$ws.Range("J74").Value = 82396.30499999999
$ws.Range("K74").Value = 13376108
$ws.Range("L74").Value = 82396.30499999999
$ws.Range("M74").Value = -13375234
$ws.Range("N74").Value = -84144.30499999999

$ws.Range("H77").Value = 8828259
$ws.Range("I77").Value = 13376108
$ws.Range("J77").Value = 82396.30499999999
$ws.Range("K77").Value = 66880540
$ws.Range("L77").Value = 411981.525
$ws.Range("M77").Value = -66876172
$ws.Range("N77").Value = -420717.525

$ws.Range("H116").Value = 2791
$ws.Range("I116").Value = 2682.8333
$ws.Range("K116").Value = 2682.8333
$ws.Range("M116").Value = -388.8332999999998

$ws.Range("H132").Value = 13572176
$ws.Range("I132").Value = 19271280
$ws.Range("K132").Value = 57813840
$ws.Range("M132").Value = -57811310

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2791
$ws.Range("I3").Value = 2682.8333
$ws.Range("K3").Value = 2682.8333
$ws.Range("M3").Value = -2568.8333

$ws.Range("H20").Value = 1955.3334
$ws.Range("I20").Value = 1805.4286
$ws.Range("K20").Value = 1805.4286
$ws.Range("M20").Value = -1558.4286

$ws.Range("H51").Value = 28000
$ws.Range("J51").Value = 28000
$ws.Range("L51").Value = 28000
$ws.Range("N51").Value = -28982

$ws.Range("H107").Value = 5116.6665
$ws.Range("I107").Value = 4940
$ws.Range("J107").Value = 6000
$ws.Range("K107").Value = 4940
$ws.Range("L107").Value = 6000
$ws.Range("M107").Value = -3020
$ws.Range("N107").Value = -9840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 15153392
$ws.Range("I58").Value = 25642744
$ws.Range("J58").Value = 2107.5557
$ws.Range("K58").Value = 25642744
$ws.Range("L58").Value = 2107.5557
$ws.Range("M58").Value = -25642541
$ws.Range("N58").Value = -2513.5557

$ws.Range("H132").Value = 22396.541
$ws.Range("I132").Value = 1257.4193
$ws.Range("K132").Value = 3772.2579
$ws.Range("M132").Value = -1242.2579

$ws.Range("H134").Value = 37942.516
$ws.Range("I134").Value = 696.9
$ws.Range("J134").Value = 120710.555
$ws.Range("K134").Value = 2090.7
$ws.Range("L134").Value = 362131.665
$ws.Range("M134").Value = 444.3000000000002
$ws.Range("N134").Value = -367201.665

$ws.Range("H136").Value = 15153392
$ws.Range("I136").Value = 25642744
$ws.Range("J136").Value = 2107.5557
$ws.Range("K136").Value = 76928232
$ws.Range("L136").Value = 6322.6671
$ws.Range("M136").Value = -76925682
$ws.Range("N136").Value = -11422.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 4472.9287
$ws.Range("I64").Value = 2312
$ws.Range("J64").Value = 4639.154
$ws.Range("K64").Value = 6936
$ws.Range("L64").Value = 13917.462
$ws.Range("M64").Value = -6666
$ws.Range("N64").Value = -14457.462

$ws.Range("H67").Value = 4472.9287
$ws.Range("I67").Value = 2312
$ws.Range("J67").Value = 4639.154
$ws.Range("K67").Value = 6936
$ws.Range("L67").Value = 13917.462
$ws.Range("M67").Value = -6000
$ws.Range("N67").Value = -15789.462

$ws.Range("H131").Value = 950.375
$ws.Range("I131").Value = 752.25
$ws.Range("J131").Value = 990
$ws.Range("K131").Value = 2256.75
$ws.Range("L131").Value = 2970
$ws.Range("M131").Value = 2783.25
$ws.Range("N131").Value = -13050

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2135.087
$ws.Range("I113").Value = 1405.7693
$ws.Range("K113").Value = 1405.7693
$ws.Range("M113").Value = 764.2307000000001

$ws.Range("H126").Value = 2703.5
$ws.Range("I126").Value = 2400
$ws.Range("J126").Value = 3007
$ws.Range("K126").Value = 7200
$ws.Range("L126").Value = 9021
$ws.Range("M126").Value = -4730
$ws.Range("N126").Value = -13961

$ws.Range("H132").Value = 56887.11
$ws.Range("I132").Value = 45715.305
$ws.Range("J132").Value = 75240.78999999999
$ws.Range("K132").Value = 137145.915
$ws.Range("L132").Value = 225722.37
$ws.Range("M132").Value = -134615.915
$ws.Range("N132").Value = -230782.37

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1046111
$ws.Range("I46").Value = 2525751
$ws.Range("K46").Value = 2525751
$ws.Range("M46").Value = -2525563

$ws.Range("H132").Value = 21976.98
$ws.Range("I132").Value = 1798.6154
$ws.Range("J132").Value = 93518.45
$ws.Range("K132").Value = 5395.8462
$ws.Range("L132").Value = 280555.35
$ws.Range("M132").Value = -2865.8462
$ws.Range("N132").Value = -285615.35

$ws.Range("H136").Value = 38756.16
$ws.Range("I136").Value = 21487.871
$ws.Range("J136").Value = 505000
$ws.Range("K136").Value = 64463.613
$ws.Range("L136").Value = 1515000
$ws.Range("M136").Value = -61913.613
$ws.Range("N136").Value = -1520100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5231.4
$ws.Range("I96").Value = 1801.5
$ws.Range("J96").Value = 7518
$ws.Range("K96").Value = 1801.5
$ws.Range("L96").Value = 7518
$ws.Range("M96").Value = -428.5
$ws.Range("N96").Value = -10264

$ws.Range("H122").Value = 2863.3333
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = $null

$ws.Range("H126").Value = 945.94116
$ws.Range("I126").Value = 576.7778
$ws.Range("J126").Value = 1361.25
$ws.Range("K126").Value = 1730.3334
$ws.Range("L126").Value = 4083.75
$ws.Range("M126").Value = 739.6666
$ws.Range("N126").Value = -9023.75

$ws.Range("H132").Value = 40455.254
$ws.Range("I132").Value = 29492.715
$ws.Range("J132").Value = 64435.812
$ws.Range("K132").Value = 88478.145
$ws.Range("L132").Value = 193307.436
$ws.Range("M132").Value = -85948.145
$ws.Range("N132").Value = -198367.436
